{"js": "const replacements = [\n  [\"61\u00f72=\", \"87\u00f72=\"],\n  [\"66\u00f73=\", \"30\u00f73=\"],\n  [\"97\u00f77=\", \"42\u00f79=\"],\n  [\"49\u00f74=\", \"91\u00f75=\"],\n  [\"76\u00f77=\", \"29\u00f76=\"],\n  [\"48\u00f72=\", \"16\u00f79=\"],\n  [\"60\u00f75=\", \"93\u00f78=\"],\n  [\"41\u00f79=\", \"14\u00f75=\"],\n  [\"38\u00f76=\", \"99\u00f77=\"],\n  [\"24\u00f78=\", \"98\u00f72=\"],\n  [\"14\u00f72=\", \"22\u00f72=\"],\n  [\"95\u00f76=\", \"51\u00f76=\"],\n  [\"27\u00f79=\", \"42\u00f76=\"],\n  [\"34\u00f76=\", \"89\u00f78=\"],\n  [\"76\u00f73=\", \"33\u00f73=\"],\n  [\"81\u00f73=\", \"81\u00f79=\"],\n  [\"12\u00f76=\", \"44\u00f77=\"],\n  [\"18\u00f73=\", \"32\u00f78=\"],\n  [\"17\u00f73=\", \"84\u00f73=\"],\n  [\"75\u00f72=\", \"83\u00f76=\"],\n  [\"12\u00f77=\", \"43\u00f73=\"],\n  [\"58\u00f72=\", \"69\u00f73=\"],\n  [\"58\u00f76=\", \"70\u00f74=\"],\n  [\"70\u00f76=\", \"81\u00f79=\"],\n  [\"45\u00f76=\", \"67\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old=\"61\u00f72=\"; new=\"87\u00f72=\"},\n    @{old=\"66\u00f73=\"; new=\"30\u00f73=\"},\n    @{old=\"97\u00f77=\"; new=\"42\u00f79=\"},\n    @{old=\"49\u00f74=\"; new=\"91\u00f75=\"},\n    @{old=\"76\u00f77=\"; new=\"29\u00f76=\"},\n    @{old=\"48\u00f72=\"; new=\"16\u00f79=\"},\n    @{old=\"60\u00f75=\"; new=\"93\u00f78=\"},\n    @{old=\"41\u00f79=\"; new=\"14\u00f75=\"},\n    @{old=\"38\u00f76=\"; new=\"99\u00f77=\"},\n    @{old=\"24\u00f78=\"; new=\"98\u00f72=\"},\n    @{old=\"14\u00f72=\"; new=\"22\u00f72=\"},\n    @{old=\"95\u00f76=\"; new=\"51\u00f76=\"},\n    @{old=\"27\u00f79=\"; new=\"42\u00f76=\"},\n    @{old=\"34\u00f76=\"; new=\"89\u00f78=\"},\n    @{old=\"76\u00f73=\"; new=\"33\u00f73=\"},\n    @{old=\"81\u00f73=\"; new=\"81\u00f79=\"},\n    @{old=\"12\u00f76=\"; new=\"44\u00f77=\"},\n    @{old=\"18\u00f73=\"; new=\"32\u00f78=\"},\n    @{old=\"17\u00f73=\"; new=\"84\u00f73=\"},\n    @{old=\"75\u00f72=\"; new=\"83\u00f76=\"},\n    @{old=\"12\u00f77=\"; new=\"43\u00f73=\"},\n    @{old=\"58\u00f72=\"; new=\"69\u00f73=\"},\n    @{old=\"58\u00f76=\"; new=\"70\u00f74=\"},\n    @{old=\"70\u00f76=\"; new=\"81\u00f79=\"},\n    @{old=\"45\u00f76=\"; new=\"67\u00f76=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
